$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''298.59'
$ws.Range("E2").Value = '''-2.61%'
$ws.Range("D3").Value = '''31.73'
$ws.Range("E3").Value = '''-1.47%'
$ws.Range("D4").Value = '''5.105'
$ws.Range("E4").Value = '''-4.00%'
$ws.Range("E5").Value = '''1.69%'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").Value = '''1.792'
$ws.Range("E6").Value = '''17.14%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '''7.743'
$ws.Range("E7").Value = '''-0.17%'
$ws.Range("D8").Value = '''3.791'
$ws.Range("E8").Value = '''2.66%'
$ws.Range("D9").Value = '''0.9260'
$ws.Range("E9").Value = '''1.87%'
$ws.Range("D10").Value = '''0.1712'
$ws.Range("E10").Value = '''2.53%'
$ws.Range("D11").Value = '''0.07276'
$ws.Range("D12").Value = '''0.07953'
$ws.Range("E12").Value = '''-1.29%'
$ws.Range("E13").Value = '''1.20%'
$ws.Range("E14").Value = '''0.46%'
$ws.Range("D15").Value = '''0.001488'
$ws.Range("E15").Value = '''-2.17%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006540'
$ws.Range("E16").Value = '''3.46%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.449'
$ws.Range("E17").Value = '''-1.32%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.218'
$ws.Range("E18").Value = '''-0.96%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3280'
$ws.Range("E19").Value = '''0.49%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = '''0.1325'
$ws.Range("E20").Value = '''-0.54%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = '''4.555'
$ws.Range("E21").Value = '''8.07%'
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").Value = '''0.04647'
$ws.Range("E22").Value = '''1.89%'
$ws.Range("E23").Value = '''-4.94%'
$ws.Range("E24").Value = '''0.06%'
$ws.Range("D25").Value = '''0.004421'
$ws.Range("E25").Value = '''-1.90%'
$ws.Range("D26").Value = '''0.0001398'
$ws.Range("E26").Value = '''19.62%'
$ws.Range("D27").Value = '''0.0001857'
$ws.Range("E27").Value = '''7.02%'
$ws.Range("E39").Value = '''1.42%'
$ws.Range("D40").Value = '''0.04540'
$ws.Range("E40").Value = '''0.92%'
$ws.Range("D41").Value = '''0.007058'
$ws.Range("E41").Value = '''-3.35%'
$ws.Range("D42").Value = '''0.1326'
$ws.Range("E42").Value = '''-2.48%'
$ws.Range("E43").Value = '''-8.88%'
$ws.Range("D44").Value = '''0.01276'
$ws.Range("E44").Value = '''-6.57%'
$ws.Range("D45").Value = '''0.00006020'
$ws.Range("E45").Value = '''-1.69%'
$ws.Range("D46").Value = '''0.7116'
$ws.Range("E46").Value = '''-62.40%'
$ws.Range("E47").Value = '''0.05%'
